$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040354665247028
$ws.Cells.Item(2, 4).Value = 1.063940348911443
$ws.Cells.Item(2, 5).Value = 1.038741983315371
$ws.Cells.Item(2, 6).Value = 1.068807601400088
$ws.Cells.Item(2, 9).Value = 1.050936796304074
$ws.Cells.Item(2, 10).Value = 1.045441722664561
$ws.Cells.Item(2, 11).Value = 1.066657630463318
$ws.Cells.Item(2, 12).Value = 1.041529151243364
$ws.Cells.Item(2, 13).Value = 1.071511794765378
$ws.Cells.Item(2, 14).Value = 1.018969147809317
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.041445113658833
$ws.Cells.Item(3, 4).Value = 1.064755127495261
$ws.Cells.Item(3, 5).Value = 1.039674764694412
$ws.Cells.Item(3, 6).Value = 1.069775293572723
$ws.Cells.Item(3, 9).Value = 1.051287187557872
$ws.Cells.Item(3, 10).Value = 1.046177050605059
$ws.Cells.Item(3, 11).Value = 1.067287518250532
$ws.Cells.Item(3, 12).Value = 1.042271684398492
$ws.Cells.Item(3, 13).Value = 1.072295159428767
$ws.Cells.Item(3, 14).Value = 1.019221291490035
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042150814551124
$ws.Cells.Item(4, 4).Value = 1.065282326312057
$ws.Cells.Item(4, 5).Value = 1.040278731125987
$ws.Cells.Item(4, 6).Value = 1.070401787075296
$ws.Cells.Item(4, 9).Value = 1.051512665462891
$ws.Cells.Item(4, 10).Value = 1.046652403729566
$ws.Cells.Item(4, 11).Value = 1.067694426764042
$ws.Cells.Item(4, 12).Value = 1.042751923291369
$ws.Cells.Item(4, 13).Value = 1.072801757019421
$ws.Cells.Item(4, 14).Value = 1.019384102841043
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042447517438207
$ws.Cells.Item(5, 4).Value = 1.065503955805186
$ws.Cells.Item(5, 5).Value = 1.040532732776218
$ws.Cells.Item(5, 6).Value = 1.070665243819074
$ws.Cells.Item(5, 9).Value = 1.05160715725143
$ws.Cells.Item(5, 10).Value = 1.046852133616438
$ws.Cells.Item(5, 11).Value = 1.067865329993585
$ws.Cells.Item(5, 12).Value = 1.042953760577163
$ws.Cells.Item(5, 13).Value = 1.073014660332251
$ws.Cells.Item(5, 14).Value = 1.019452466610427
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042497336708838
$ws.Cells.Item(6, 4).Value = 1.06554116806459
$ws.Cells.Item(6, 5).Value = 1.040575386284494
$ws.Cells.Item(6, 6).Value = 1.070709483973585
$ws.Cells.Item(6, 9).Value = 1.051623005292331
$ws.Cells.Item(6, 10).Value = 1.046885662775081
$ws.Cells.Item(6, 11).Value = 1.067894015936647
$ws.Cells.Item(6, 12).Value = 1.042987646707496
$ws.Cells.Item(6, 13).Value = 1.073050403598273
$ws.Cells.Item(6, 14).Value = 1.019463940368525
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04215477900602
$ws.Cells.Item(7, 4).Value = 1.065285287756311
$ws.Cells.Item(7, 5).Value = 1.040282124738826
$ws.Cells.Item(7, 6).Value = 1.07040530708802
$ws.Cells.Item(7, 9).Value = 1.051513929242376
$ws.Cells.Item(7, 10).Value = 1.04665507295676
$ws.Cells.Item(7, 11).Value = 1.067696711014999
$ws.Cells.Item(7, 12).Value = 1.042754620467929
$ws.Cells.Item(7, 13).Value = 1.072804602120787
$ws.Cells.Item(7, 14).Value = 1.019385016643927
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040723165308154
$ws.Cells.Item(8, 4).Value = 1.064215710320183
$ws.Cells.Item(8, 5).Value = 1.039057139294304
$ws.Cells.Item(8, 6).Value = 1.069134568407921
$ws.Cells.Item(8, 9).Value = 1.051055471099299
$ws.Cells.Item(8, 10).Value = 1.045690323988888
$ws.Cells.Item(8, 11).Value = 1.066870642643825
$ws.Cells.Item(8, 12).Value = 1.041780141271349
$ws.Cells.Item(8, 13).Value = 1.071776596640754
$ws.Cells.Item(8, 14).Value = 1.019054431776733
$ws.Cells.Item(9, 2).Value = 1.019999999999999
$ws.Cells.Item(9, 3).Value = 1.038201286374832
$ws.Cells.Item(9, 4).Value = 1.062330876398452
$ws.Cells.Item(9, 5).Value = 1.036901592000276
$ws.Cells.Item(9, 6).Value = 1.066897932863454
$ws.Cells.Item(9, 9).Value = 1.050238054152035
$ws.Cells.Item(9, 10).Value = 1.043986844360197
$ws.Cells.Item(9, 11).Value = 1.065409885974595
$ws.Cells.Item(9, 12).Value = 1.040061231471274
$ws.Cells.Item(9, 13).Value = 1.069962905920067
$ws.Cells.Item(9, 14).Value = 1.018469281906253
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.036520547662603
$ws.Cells.Item(10, 4).Value = 1.061074292617718
$ws.Cells.Item(10, 5).Value = 1.035466614015469
$ws.Cells.Item(10, 6).Value = 1.065408602061996
$ws.Cells.Item(10, 9).Value = 1.049686694878113
$ws.Cells.Item(10, 10).Value = 1.042848860092088
$ws.Cells.Item(10, 11).Value = 1.064432632496464
$ws.Cells.Item(10, 12).Value = 1.038914119398189
$ws.Cells.Item(10, 13).Value = 1.068752316064401
$ws.Cells.Item(10, 14).Value = 1.018077429497981
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.035792883846964
$ws.Cells.Item(11, 4).Value = 1.060530179330679
$ws.Cells.Item(11, 5).Value = 1.034845741935686
$ws.Cells.Item(11, 6).Value = 1.064764128897806
$ws.Cells.Item(11, 9).Value = 1.049446430350058
$ws.Cells.Item(11, 10).Value = 1.042355546842398
$ws.Cells.Item(11, 11).Value = 1.064008665654524
$ws.Cells.Item(11, 12).Value = 1.038417129798814
$ws.Cells.Item(11, 13).Value = 1.068227775094442
$ws.Cells.Item(11, 14).Value = 1.017907338655451
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.035522612468815
$ws.Cells.Item(12, 4).Value = 1.060328071364518
$ws.Cells.Item(12, 5).Value = 1.034615195026577
$ws.Cells.Item(12, 6).Value = 1.064524805915517
$ws.Cells.Item(12, 9).Value = 1.049356956753864
$ws.Cells.Item(12, 10).Value = 1.042172224274353
$ws.Cells.Item(12, 11).Value = 1.063851064167723
$ws.Cells.Item(12, 12).Value = 1.038232483249708
$ws.Cells.Item(12, 13).Value = 1.068032885300228
$ws.Cells.Item(12, 14).Value = 1.01784409688492
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.035580585913463
$ws.Cells.Item(13, 4).Value = 1.060371424203142
$ws.Cells.Item(13, 5).Value = 1.034664644833406
$ws.Cells.Item(13, 6).Value = 1.064576138648052
$ws.Cells.Item(13, 9).Value = 1.049376159490022
$ws.Cells.Item(13, 10).Value = 1.042211551393832
$ws.Cells.Item(13, 11).Value = 1.063884875706783
$ws.Cells.Item(13, 12).Value = 1.038272092484395
$ws.Cells.Item(13, 13).Value = 1.068074692167516
$ws.Cells.Item(13, 14).Value = 1.017857665284691
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.03577054282138
$ws.Cells.Item(14, 4).Value = 1.060513473021182
$ws.Cells.Item(14, 5).Value = 1.034826683374241
$ws.Cells.Item(14, 6).Value = 1.064744345096932
$ws.Cells.Item(14, 9).Value = 1.049439039098782
$ws.Cells.Item(14, 10).Value = 1.04234039505742
$ws.Cells.Item(14, 11).Value = 1.063995640739627
$ws.Cells.Item(14, 12).Value = 1.038401867724343
$ws.Cells.Item(14, 13).Value = 1.06821166650075
$ws.Cells.Item(14, 14).Value = 1.01790211234472
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.035887583555473
$ws.Cells.Item(15, 4).Value = 1.060600993994599
$ws.Cells.Item(15, 5).Value = 1.034926530293263
$ws.Cells.Item(15, 6).Value = 1.064847991008118
$ws.Cells.Item(15, 9).Value = 1.049477751004127
$ws.Cells.Item(15, 10).Value = 1.042419768747516
$ws.Cells.Item(15, 11).Value = 1.064063870664785
$ws.Cells.Item(15, 12).Value = 1.038481820902138
$ws.Cells.Item(15, 13).Value = 1.06829605403677
$ws.Cells.Item(15, 14).Value = 1.01792948937098
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.036568842508075
$ws.Cells.Item(16, 4).Value = 1.061110403584325
$ws.Cells.Item(16, 5).Value = 1.035507829430824
$ws.Cells.Item(16, 6).Value = 1.065451382442134
$ws.Cells.Item(16, 9).Value = 1.049702608374957
$ws.Cells.Item(16, 10).Value = 1.042881587887217
$ws.Cells.Item(16, 11).Value = 1.064460752748168
$ws.Cells.Item(16, 12).Value = 1.038947097005608
$ws.Cells.Item(16, 13).Value = 1.068787120823353
$ws.Cells.Item(16, 14).Value = 1.018088709104513
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.036996206415244
$ws.Cells.Item(17, 4).Value = 1.06142994203462
$ws.Cells.Item(17, 5).Value = 1.035872592354389
$ws.Cells.Item(17, 6).Value = 1.065829985854265
$ws.Cells.Item(17, 9).Value = 1.049843247935164
$ws.Cells.Item(17, 10).Value = 1.043171125188716
$ws.Cells.Item(17, 11).Value = 1.064709489852801
$ws.Cells.Item(17, 12).Value = 1.039238876750139
$ws.Cells.Item(17, 13).Value = 1.069095061198717
$ws.Cells.Item(17, 14).Value = 1.018188472064164
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037245490938189
$ws.Cells.Item(18, 4).Value = 1.061616322961514
$ws.Cells.Item(18, 5).Value = 1.036085398979608
$ws.Cells.Item(18, 6).Value = 1.066050859016454
$ws.Cells.Item(18, 9).Value = 1.049925133724958
$ws.Cells.Item(18, 10).Value = 1.043339953515662
$ws.Cells.Item(18, 11).Value = 1.064854495858356
$ws.Cells.Item(18, 12).Value = 1.039409039683938
$ws.Cells.Item(18, 13).Value = 1.069274644054054
$ws.Cells.Item(18, 14).Value = 1.018246621952473
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037330492302501
$ws.Cells.Item(19, 4).Value = 1.06167987394825
$ws.Cells.Item(19, 5).Value = 1.036157968438224
$ws.Cells.Item(19, 6).Value = 1.066126177862525
$ws.Cells.Item(19, 9).Value = 1.049953029746921
$ws.Cells.Item(19, 10).Value = 1.043397510477306
$ws.Cells.Item(19, 11).Value = 1.064903925895615
$ws.Cells.Item(19, 12).Value = 1.039467056225349
$ws.Cells.Item(19, 13).Value = 1.069335871466244
$ws.Cells.Item(19, 14).Value = 1.01826644274814
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.036950353243377
$ws.Cells.Item(20, 4).Value = 1.061395658624828
$ws.Cells.Item(20, 5).Value = 1.035833451925762
$ws.Cells.Item(20, 6).Value = 1.065789361130345
$ws.Cells.Item(20, 9).Value = 1.049828173834003
$ws.Cells.Item(20, 10).Value = 1.043140066150687
$ws.Cells.Item(20, 11).Value = 1.064682810795468
$ws.Cells.Item(20, 12).Value = 1.039207574364206
$ws.Cells.Item(20, 13).Value = 1.069062025579353
$ws.Cells.Item(20, 14).Value = 1.018177772595987
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.035714604828355
$ws.Cells.Item(21, 4).Value = 1.060471643181968
$ws.Cells.Item(21, 5).Value = 1.034778965057185
$ws.Cells.Item(21, 6).Value = 1.064694810746491
$ws.Cells.Item(21, 9).Value = 1.049420528932685
$ws.Cells.Item(21, 10).Value = 1.042302456131535
$ws.Cells.Item(21, 11).Value = 1.06396302654483
$ws.Cells.Item(21, 12).Value = 1.038363653323813
$ws.Cells.Item(21, 13).Value = 1.068171332404552
$ws.Cells.Item(21, 14).Value = 1.017889025516982
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03493772927355
$ws.Cells.Item(22, 4).Value = 1.059890677905143
$ws.Cells.Item(22, 5).Value = 1.034116387541168
$ws.Cells.Item(22, 6).Value = 1.064006988446064
$ws.Cells.Item(22, 9).Value = 1.049162903520314
$ws.Cells.Item(22, 10).Value = 1.041775331037909
$ws.Cells.Item(22, 11).Value = 1.063509767827179
$ws.Cells.Item(22, 12).Value = 1.037832800855363
$ws.Cells.Item(22, 13).Value = 1.067611017991525
$ws.Cells.Item(22, 14).Value = 1.017707117436309
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035349557557623
$ws.Cells.Item(23, 4).Value = 1.06019865840997
$ws.Cells.Item(23, 5).Value = 1.034467592558305
$ws.Cells.Item(23, 6).Value = 1.064371581329526
$ws.Cells.Item(23, 9).Value = 1.04929960099977
$ws.Cells.Item(23, 10).Value = 1.042054816179979
$ws.Cells.Item(23, 11).Value = 1.063750115266187
$ws.Cells.Item(23, 12).Value = 1.038114239106698
$ws.Cells.Item(23, 13).Value = 1.067908079625915
$ws.Cells.Item(23, 14).Value = 1.017803584579756
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.036971072285649
$ws.Cells.Item(24, 4).Value = 1.061411149823357
$ws.Cells.Item(24, 5).Value = 1.035851137656636
$ws.Cells.Item(24, 6).Value = 1.065807717570153
$ws.Cells.Item(24, 9).Value = 1.049834985625248
$ws.Cells.Item(24, 10).Value = 1.043154100559398
$ws.Cells.Item(24, 11).Value = 1.064694866154477
$ws.Cells.Item(24, 12).Value = 1.039221718649107
$ws.Cells.Item(24, 13).Value = 1.069076953057357
$ws.Cells.Item(24, 14).Value = 1.018182607349155
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038853159843231
$ws.Cells.Item(25, 4).Value = 1.062818159019021
$ws.Cells.Item(25, 5).Value = 1.037458491714048
$ws.Cells.Item(25, 6).Value = 1.067475849021882
$ws.Cells.Item(25, 9).Value = 1.050450507363125
$ws.Cells.Item(25, 10).Value = 1.044427645636857
$ws.Cells.Item(25, 11).Value = 1.065788130513794
$ws.Cells.Item(25, 12).Value = 1.040505817935361
$ws.Cells.Item(25, 13).Value = 1.070432048134732
$ws.Cells.Item(25, 14).Value = 1.018620866583206
